$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

Set-TextCell "D2" "23.489.96"
Set-TextCell "E2" "  -2.19%  "
Set-TextCell "D3" "1.612.88"
Set-TextCell "E3" "  -2.78%  "
Set-TextCell "D4" "1.008"
Set-TextCell "E4" "  +0.71%  "
Set-TextCell "D5" "1.005"
Set-TextCell "E5" "  +0.37%  "
Set-TextCell "D6" "305.78"
Set-TextCell "E6" "  -1.40%  "
Set-TextCell "D7" "0.3830"
Set-TextCell "E7" "  -1.85%  "
Set-TextCell "D8" "0.3748"
Set-TextCell "E8" "  -2.94%  "
Set-TextCell "D9" "49.46"
Set-TextCell "E9" "  -3.54%  "
Set-TextCell "D10" "1.008"
Set-TextCell "E10" "  +0.71%  "
Set-TextCell "D11" "1.296"
Set-TextCell "E11" "  -4.97%  "
Set-TextCell "D12" "0.08223"
Set-TextCell "E12" "  -3.24%  "
Set-TextCell "D13" "23.16"
Set-TextCell "E13" "  -3.09%  "
Set-TextCell "D14" "6.781"
Set-TextCell "E14" "  -5.97%  "
Set-TextCell "D15" "7.629"
Set-TextCell "E15" "  -4.79%  "
Set-TextCell "D16" "0.00001277"
Set-TextCell "E16" "  -2.70%  "
Set-TextCell "D17" "1.617.71"
Set-TextCell "E17" "  -2.36%  "
Set-TextCell "D18" "92.63"
Set-TextCell "E18" "  -2.03%  "
Set-TextCell "D19" "0.06897"
Set-TextCell "E19" "  -1.27%  "
Set-TextCell "D20" "18.93"
Set-TextCell "E20" "  -4.84%  "
Set-TextCell "D21" "6.734"
Set-TextCell "E21" "  -3.36%  "
Set-TextCell "D22" "1.006"
Set-TextCell "E22" "  +0.57%  "
Set-TextCell "D23" "13.32"
Set-TextCell "E23" "  -2.40%  "
Set-TextCell "D24" "23.495.32"
Set-TextCell "E24" "  -2.14%  "
Set-TextCell "D25" "2.404"
Set-TextCell "E25" "  -3.37%  "
Set-TextCell "D26" "2.856"
Set-TextCell "E26" "  -8.17%  "
Set-TextCell "D27" "21.50"
Set-TextCell "E27" "  -3.11%  "
Set-TextCell "D28" "151.92"
Set-TextCell "E28" "  -1.06%  "
Set-TextCell "D29" "5.394"
Set-TextCell "E29" "  +1.40%  "
Set-TextCell "D30" "7.872"
Set-TextCell "E30" "  -0.91%  "
Set-TextCell "D31" "134.41"
Set-TextCell "E31" "  -3.80%  "
Set-TextCell "D32" "2.453"
Set-TextCell "E32" "  -1.37%  "
Set-TextCell "D33" "1.805.20"
Set-TextCell "E33" "  -1.83%  "
Set-TextCell "D34" "0.9643"
Set-TextCell "E34" "  -7.16%  "
Set-TextCell "D35" "0.07700"
Set-TextCell "E35" "  -5.34%  "
Set-TextCell "D36" "0.02841"
Set-TextCell "E36" "  -5.10%  "
Set-TextCell "D37" "6.472"
Set-TextCell "E37" "  -3.56%  "
Set-TextCell "D38" "0.2605"
Set-TextCell "E38" "  -3.44%  "
Set-TextCell "D39" "10.21"
Set-TextCell "E39" "  -7.85%  "
Set-TextCell "D40" "0.08966"
Set-TextCell "E40" "  -2.01%  "
Set-TextCell "D41" "0.7323"
Set-TextCell "E41" "  -2.98%  "
Set-TextCell "D42" "13.09"
Set-TextCell "E42" "  -4.11%  "
Set-TextCell "D43" "1.395"
Set-TextCell "E43" "  -1.77%  "
Set-TextCell "D44" "15.98"
Set-TextCell "E44" "  -2.25%  "
Set-TextCell "D45" "0.6746"
Set-TextCell "E45" "  -3.76%  "
Set-TextCell "D46" "2.358"
Set-TextCell "E46" "  -5.17%  "
Set-TextCell "B47" "Frax"
Set-TextCell "C47" "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextCell "D47" "1.003"
Set-TextCell "E47" "  +0.25%  "
Set-TextCell "B48" "PancakeSwap"
Set-TextCell "C48" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell "D48" "4.031"
Set-TextCell "E48" "  -1.47%  "
Set-TextCell "D49" "0.08092"
Set-TextCell "E49" "  -2.27%  "
Set-TextCell "D50" "132.27"
Set-TextCell "E50" "  -2.34%  "
Set-TextCell "D51" "1.190"
Set-TextCell "E51" "  -3.57%  "
